# Apply updated "dSF" (column F) values as part of a repull/recalculation of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value  = -4
$ws.Range("F3").Value  = 5
$ws.Range("F5").Value  = 7
$ws.Range("F8").Value  = -3
$ws.Range("F9").Value  = 0
$ws.Range("F11").Value = -4
$ws.Range("F15").Value = -3
$ws.Range("F17").Value = -8
